# Applies the "script update" diff to the betexplorer primera-rfef-group-1
# 2023-2024 sheet:
#   - swap the match-detail columns (F:V) between several pairs of rows
#     (the match order within a kickoff-day block changed, A:E - the
#     index/meta columns - stay put)
#   - rotate three rows (19,20,21) that got reordered as a 3-cycle
#   - append one new match row (116) with full odds/url data
#
# Columns F..V = 6..22 (home name .. match url); columns A..E (index,
# pais, torneio, temporada, data_partida) are left untouched for the
# swapped/rotated rows because the diff never touches them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colStart = 6   # F
$colEnd   = 22  # V

function Swap-RowRange($rowA, $rowB, $cStart, $cEnd) {
    for ($c = $cStart; $c -le $cEnd; $c++) {
        $cellA = $ws.Cells.Item($rowA, $c)
        $cellB = $ws.Cells.Item($rowB, $c)
        $valA = $cellA.Value()
        $valB = $cellB.Value()
        $cellA.Value = $valB
        $cellB.Value = $valA
    }
}

# ---- simple two-row swaps (F:V only) ----
$swapPairs = @(
    @(2, 3),
    @(12, 13),
    @(14, 15),
    @(30, 31),
    @(34, 35),
    @(36, 37),
    @(56, 57),
    @(72, 73),
    @(80, 81),
    @(92, 93)
)

foreach ($pair in $swapPairs) {
    Swap-RowRange $pair[0] $pair[1] $colStart $colEnd
}

# ---- 3-way rotation for rows 19, 20, 21 (F:V only) ----
# new19 = old21, new20 = old19, new21 = old20
$row19 = @{}
$row20 = @{}
$row21 = @{}
for ($c = $colStart; $c -le $colEnd; $c++) {
    $row19[$c] = $ws.Cells.Item(19, $c).Value()
    $row20[$c] = $ws.Cells.Item(20, $c).Value()
    $row21[$c] = $ws.Cells.Item(21, $c).Value()
}
for ($c = $colStart; $c -le $colEnd; $c++) {
    $ws.Cells.Item(19, $c).Value = $row21[$c]
    $ws.Cells.Item(20, $c).Value = $row19[$c]
    $ws.Cells.Item(21, $c).Value = $row20[$c]
}

# ---- append new row 116 ----
$newRow = 116

$ws.Cells.Item($newRow, 1).Value = 115
$ws.Cells.Item($newRow, 2).Value = "spain"
$ws.Cells.Item($newRow, 3).Value = "primera-rfef-group-1"
$ws.Cells.Item($newRow, 4).Value = "2023-2024"
$ws.Cells.Item($newRow, 5).Value = 45242.5
$ws.Cells.Item($newRow, 6).Value = "Tarazona"
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = "Teruel"
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = 2.14
$ws.Cells.Item($newRow, 11).Value = "09/11/2023 09:13"
$ws.Cells.Item($newRow, 12).Value = 2.68
$ws.Cells.Item($newRow, 13).Value = "12/11/2023 11:53"
$ws.Cells.Item($newRow, 14).Value = 2.96
$ws.Cells.Item($newRow, 15).Value = "09/11/2023 09:13"
$ws.Cells.Item($newRow, 16).Value = 2.63
$ws.Cells.Item($newRow, 17).Value = "12/11/2023 11:53"
$ws.Cells.Item($newRow, 18).Value = 3.36
$ws.Cells.Item($newRow, 19).Value = "09/11/2023 09:13"
$ws.Cells.Item($newRow, 20).Value = 3.28
$ws.Cells.Item($newRow, 21).Value = "12/11/2023 11:53"
$ws.Cells.Item($newRow, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-1/tarazona-teruel/pSIXUlD8/"

# Match existing formatting: column A uses the bold/bordered style, column E
# uses the custom date-time number format - copy both from the row above.
$ws.Cells.Item(115, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122)
$ws.Cells.Item(115, 5).Copy()
$ws.Cells.Item($newRow, 5).PasteSpecial(-4122)

$excel.CutCopyMode = 0
